$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new row at row 9 (pushes existing rows 9+ down by one,
#    carrying the row-9 formatting onto the new blank row).
$ws.Rows("9").Insert()

# 2. Populate the new row 9 with the MAX 10 FPGA soldering tip (point "2."),
#    and renumber the two tips that were pushed down ("2." -> "3.", "3." -> "4.").
$ws.Range("A9").Value = "2. Soldering the MAX 10 FPGA - Place a small amount of solder paste on the large central pad and then place the chip accordingly with gel flux applied to the pads as per point 3. Drag solder the pins first and then heat the chip from the other side of the board with a hot air gun to hopefully melt the solder paste. I'm not sure this is the correct technique but thats the process I used."
$ws.Range("A10").Value = "3. I recommend using a gel like flux for the SMD chips. Run plenty along each row of pads and then align the chip correctly, the gel helps keep the chip in place. Tag each corner by soldering one or two pins and then drag solder the rest. See guide below"
$ws.Range("A11").Value = "4. All the passives (capacitors, resistors), comparator and regulator etc are best done with solder flux and a hot air gun. I would avoid doing this to the SW1 switch however as it may melt it. The switch is best done with the soldering iron"

# 3. The row insert does not shift the worksheet's hyperlink anchors, so
#    rebuild them (same targets / same display text, new row numbers).
#    Hyperlinks.Add() resets cell formatting to the generic "Hyperlink" cell
#    style, so snapshot each cell's exact formatting into scratch cells
#    first and paste it back (PasteSpecial formats) once the links are
#    rebuilt.
$hlCells = @("G16","G17","G19","G25","G28","G29","G30","G31","A5")
$scratchRow = 101
$scratchAddr = @{}
foreach ($addr in $hlCells) {
    $scratchAddr[$addr] = "Z" + $scratchRow
    $ws.Range($addr).Copy()
    $ws.Range($scratchAddr[$addr]).PasteSpecial(-4122)
    $scratchRow = $scratchRow + 1
}

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("G16"), "https://www.digikey.co.uk/product-detail/en/on-semiconductor/BAT54HT1G/BAT54HT1GOSCT-ND/917809")
$ws.Hyperlinks.Add($ws.Range("G17"), "https://www.digikey.co.uk/product-detail/en/nexperia-usa-inc/74HCT125PW,118/1727-4087-1-ND/1965389")
$ws.Hyperlinks.Add($ws.Range("G19"), "https://www.digikey.co.uk/product-detail/en/nidec-copal-electronics/CHS-04TB/563-1008-1-ND/948417")
# G25 keeps its original "display" text override (full tracking URL), same as the source cell text.
$g25Target = "https://www.digikey.co.uk/product-detail/en/diodes-incorporated/AP2114H-3.3TRG1/AP2114H-3.3TRG1DICT-ND/4505142?utm_adgroup=PMIC%20-%20Voltage%20Regulators%20-%20Linear&utm_source=google&utm_medium=cpc&utm_campaign=Google%20Shopping_Integrated%20Circuits%20%28ICs%29&utm_term=&productid=4505142&gclid=EAIaIQobChMIkoq-teuO6gIVC4myCh0O7wooEAQYASABEgJpyfD_BwE"
$ws.Hyperlinks.Add($ws.Range("G25"), $g25Target, "", "", $g25Target)
$ws.Hyperlinks.Add($ws.Range("G28"), "https://www.digikey.co.uk/products/en?keywords=10M08SCE144C8G")
$ws.Hyperlinks.Add($ws.Range("G29"), "https://www.digikey.co.uk/product-detail/en/nexperia-usa-inc/74LVC4245APW-118/1727-4308-1-ND/2209918")
$ws.Hyperlinks.Add($ws.Range("G30"), "https://hobbycomponents.com/connectors/439-01-254mm-40way-sil-turned-pin-m-m-headers-pack-of-5")
$ws.Hyperlinks.Add($ws.Range("G31"), "https://www.digikey.co.uk/product-detail/en/microchip-technology/MIC7221YM5-TR/576-2901-1-ND")
# A5 keeps its original "display" text override too (full tracking URL).
$a5Target = "https://www.ebay.co.uk/i/292483310138?chn=ps&norover=1&mkevt=1&mkrid=710-134428-41853-0&mkcid=2&itemid=292483310138&targetid=908661474856&device=c&mktype=pla&googleloc=9046613&poi=&campaignid=10195651586&mkgroupid=107296210212&rlsatarget=pla-908661474856&abcId=1145987&merchantid=7398364&gclid=EAIaIQobChMI6ru5t-WQ6gIVzoKyCh1oBwxIEAQYCiABEgImJvD_BwE"
$ws.Hyperlinks.Add($ws.Range("A5"), $a5Target, "", "", $a5Target)

foreach ($addr in $hlCells) {
    $ws.Range($scratchAddr[$addr]).Copy()
    $ws.Range($addr).PasteSpecial(-4122)
}
$ws.Range("Z101:Z109").Clear()

# 4. Shift the ElectronULA defined name to follow the table (now one row lower).
$wb.Names.Item(1).RefersTo = "=Sheet1!`$A`$15:`$I`$32"
